$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells we touch are treated as plain text (matches original inlineStr usage)
# so Excel doesn't auto-convert numeric-looking or date-looking strings.
$ws.Range("E60").NumberFormat = "@"
$ws.Range("E60").Value = "807789682"

$textRange = $ws.Range("A61:L61")
$textRange.NumberFormat = "@"

$ws.Range("A61").Value = "6271"
$ws.Range("B61").Value = "6/26/2025"
$ws.Range("C61").Value = "ARGERICH 740"
$ws.Range("D61").Value = "7"
$ws.Range("E61").Value = "807789686"
$ws.Range("F61").Value = "Optical Power"
$ws.Range("G61").Value = "Pendiente"
$ws.Range("H61").Value = "Picada"
$ws.Range("I61").Value = "1"
$ws.Range("J61").Value = "Cambio"
$ws.Range("K61").Value = "Sin equipos"
$ws.Range("L61").Value = "Pasante"

$ws.Range("M61").Value = -58.474467
$ws.Range("N61").Value = -34.624161
